$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Plain text cells (labels) -- write directly, Excel keeps them as text
$textCells = @(
    @('A2', 'Total'),
    @('E2', 'Total Actual'),
    @('A3', 'Total'),
    @('E3', 'Actual'),
    @('A4', 'Total'),
    @('E4', 'Total plan'),
    @('A5', 'Total'),
    @('E5', 'Plan'),
    @('A7', 'MGP12A-AC891-040'),
    @('E7', 'Total Actual'),
    @('A8', 'MGP12A-AC891-040'),
    @('E8', 'Actual'),
    @('A9', 'MGP12A-AC891-040'),
    @('E9', 'Total plan'),
    @('A11', 'MGP12A-AC891-040'),
    @('E11', 'Plan'),
    @('A12', 'MGP12A-AC891-030'),
    @('E12', 'Total Actual'),
    @('A13', 'MGP12A-AC891-030'),
    @('E13', 'Actual'),
    @('A14', 'MGP12A-AC891-030'),
    @('E14', 'Total plan'),
    @('A15', 'MGP12A-AC891-030'),
    @('E15', 'Plan')
)

foreach ($item in $textCells) {
    $ws.Range($item[0]).Value = $item[1]
}

# Numeric-looking cells that must be stored as TEXT (matches source inlineStr cells)
$numericTextCells = @(
    @('B2', '104999'),
    @('G2', '0'),
    @('H2', '0'),
    @('I2', '0'),
    @('J2', '0'),
    @('K2', '0'),
    @('L2', '0'),
    @('M2', '0'),
    @('N2', '0'),
    @('O2', '0'),
    @('P2', '0'),
    @('Q2', '0'),
    @('R2', '0'),
    @('S2', '0'),
    @('T2', '0'),
    @('U2', '0'),
    @('V2', '0'),
    @('W2', '0'),
    @('X2', '0'),
    @('Y2', '0'),
    @('Z2', '0'),
    @('AA2', '0'),
    @('AB2', '0'),
    @('AC2', '0'),
    @('AD2', '104999'),
    @('AE2', '104999'),
    @('AF2', '104999'),
    @('AG2', '104999'),
    @('AH2', '104999'),
    @('AI2', '104999'),
    @('AD3', '104999'),
    @('C4', '15000'),
    @('G4', '0'),
    @('H4', '0'),
    @('I4', '0'),
    @('J4', '0'),
    @('K4', '0'),
    @('L4', '0'),
    @('M4', '0'),
    @('N4', '0'),
    @('O4', '0'),
    @('P4', '0'),
    @('Q4', '0'),
    @('R4', '0'),
    @('S4', '0'),
    @('T4', '0'),
    @('U4', '0'),
    @('V4', '0'),
    @('W4', '0'),
    @('X4', '0'),
    @('Y4', '0'),
    @('Z4', '0'),
    @('AA4', '0'),
    @('AB4', '0'),
    @('AC4', '0'),
    @('AD4', '10000'),
    @('AE4', '15000'),
    @('AF4', '15000'),
    @('AG4', '15000'),
    @('AH4', '15000'),
    @('AI4', '15000'),
    @('AD5', '10000'),
    @('AE5', '5000'),
    @('B7', '4998'),
    @('G7', '0'),
    @('H7', '0'),
    @('I7', '0'),
    @('J7', '0'),
    @('K7', '0'),
    @('L7', '0'),
    @('M7', '0'),
    @('N7', '0'),
    @('O7', '0'),
    @('P7', '0'),
    @('Q7', '0'),
    @('R7', '0'),
    @('S7', '0'),
    @('T7', '0'),
    @('U7', '0'),
    @('V7', '0'),
    @('W7', '0'),
    @('X7', '0'),
    @('Y7', '0'),
    @('Z7', '0'),
    @('AA7', '0'),
    @('AB7', '0'),
    @('AC7', '0'),
    @('AD7', '4998'),
    @('AE7', '4998'),
    @('AF7', '4998'),
    @('AG7', '4998'),
    @('AH7', '4998'),
    @('AI7', '4998'),
    @('AD8', '4998'),
    @('C9', '5000'),
    @('G9', '0'),
    @('H9', '0'),
    @('I9', '0'),
    @('J9', '0'),
    @('K9', '0'),
    @('L9', '0'),
    @('M9', '0'),
    @('N9', '0'),
    @('O9', '0'),
    @('P9', '0'),
    @('Q9', '0'),
    @('R9', '0'),
    @('S9', '0'),
    @('T9', '0'),
    @('U9', '0'),
    @('V9', '0'),
    @('W9', '0'),
    @('X9', '0'),
    @('Y9', '0'),
    @('Z9', '0'),
    @('AA9', '0'),
    @('AB9', '0'),
    @('AC9', '0'),
    @('AD9', '0'),
    @('AE9', '5000'),
    @('AF9', '5000'),
    @('AG9', '5000'),
    @('AH9', '5000'),
    @('AI9', '5000'),
    @('AE11', '5000'),
    @('B12', '100001'),
    @('G12', '0'),
    @('H12', '0'),
    @('I12', '0'),
    @('J12', '0'),
    @('K12', '0'),
    @('L12', '0'),
    @('M12', '0'),
    @('N12', '0'),
    @('O12', '0'),
    @('P12', '0'),
    @('Q12', '0'),
    @('R12', '0'),
    @('S12', '0'),
    @('T12', '0'),
    @('U12', '0'),
    @('V12', '0'),
    @('W12', '0'),
    @('X12', '0'),
    @('Y12', '0'),
    @('Z12', '0'),
    @('AA12', '0'),
    @('AB12', '0'),
    @('AC12', '0'),
    @('AD12', '100001'),
    @('AE12', '100001'),
    @('AF12', '100001'),
    @('AG12', '100001'),
    @('AH12', '100001'),
    @('AI12', '100001'),
    @('AD13', '100001'),
    @('C14', '10000'),
    @('G14', '0'),
    @('H14', '0'),
    @('I14', '0'),
    @('J14', '0'),
    @('K14', '0'),
    @('L14', '0'),
    @('M14', '0'),
    @('N14', '0'),
    @('O14', '0'),
    @('P14', '0'),
    @('Q14', '0'),
    @('R14', '0'),
    @('S14', '0'),
    @('T14', '0'),
    @('U14', '0'),
    @('V14', '0'),
    @('W14', '0'),
    @('X14', '0'),
    @('Y14', '0'),
    @('Z14', '0'),
    @('AA14', '0'),
    @('AB14', '0'),
    @('AC14', '0'),
    @('AD14', '10000'),
    @('AE14', '10000'),
    @('AF14', '10000'),
    @('AG14', '10000'),
    @('AH14', '10000'),
    @('AI14', '10000'),
    @('AD15', '10000')
)

foreach ($item in $numericTextCells) {
    $rng = $ws.Range($item[0])
    $rng.NumberFormat = "@"
    $rng.Value = $item[1]
}

